# Flip the sign of the "material recycled" figures on every year sheet.
# Each sheet has the same layout: row 2/4 hold C2/B4/C4/E4 and row 5 holds D5
# as the only cells that may carry non-zero numeric data; zero cells stay zero.
$wb = $excel.ActiveWorkbook

$targets = @("C2", "B4", "C4", "E4", "D5")

foreach ($ws in $wb.Worksheets) {
    foreach ($addr in $targets) {
        $cell = $ws.Range($addr)
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne 0) {
            $cell.Value2 = -$val
        }
    }
}
